$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 773.7143
$ws.Range("I28").Value = 796.1667
$ws.Range("J28").Value = 743.7778
$ws.Range("K28").Value = 796.1667
$ws.Range("L28").Value = 743.7778
$ws.Range("M28").Value = -311.1667
$ws.Range("N28").Value = -1713.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 423.8889
$ws.Range("I107").Value = 423.8889
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 423.8889
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1496.1111
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 693014.25
$ws.Range("I132").Value = 2630.5933
$ws.Range("J132").Value = 4087400.5
$ws.Range("K132").Value = 7891.7799
$ws.Range("L132").Value = 12262201.5
$ws.Range("M132").Value = -5361.7799
$ws.Range("N132").Value = -12267261.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 42056.4
$ws.Range("I135").Value = 43517.082
$ws.Range("K135").Value = 391653.738
$ws.Range("M135").Value = -389118.738

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2043735.9
$ws.Range("I137").Value = 2441191
$ws.Range("J137").Value = 6778.875
$ws.Range("K137").Value = 7323573
$ws.Range("L137").Value = 20336.625
$ws.Range("M137").Value = -7321023
$ws.Range("N137").Value = -25436.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2819878.8
$ws.Range("I138").Value = 2350.36
$ws.Range("J138").Value = 4351144.5
$ws.Range("K138").Value = 7051.08
$ws.Range("L138").Value = 13053433.5
$ws.Range("M138").Value = -1911.08
$ws.Range("N138").Value = -13063713.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3489.1853
$ws.Range("I141").Value = 3531.3333
$ws.Range("J141").Value = 3341.6667
$ws.Range("K141").Value = 10593.9999
$ws.Range("L141").Value = 10025.0001
$ws.Range("M141").Value = -5413.999899999999
$ws.Range("N141").Value = -20385.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 37112812
$ws.Range("I61").Value = 52685890
$ws.Range("J61").Value = 126750
$ws.Range("K61").Value = 52685890
$ws.Range("L61").Value = 126750
$ws.Range("M61").Value = -52685678
$ws.Range("N61").Value = -127174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5850738.5
$ws.Range("I74").Value = 7409364
$ws.Range("J74").Value = 5891.6665
$ws.Range("K74").Value = 7409364
$ws.Range("L74").Value = 5891.6665
$ws.Range("M74").Value = -7408490
$ws.Range("N74").Value = -7639.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5850738.5
$ws.Range("I77").Value = 7409364
$ws.Range("J77").Value = 5891.6665
$ws.Range("K77").Value = 37046820
$ws.Range("L77").Value = 29458.3325
$ws.Range("M77").Value = -37042452
$ws.Range("N77").Value = -38194.3325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 37112812
$ws.Range("I136").Value = 52685890
$ws.Range("J136").Value = 126750
$ws.Range("K136").Value = 158057670
$ws.Range("L136").Value = 380250
$ws.Range("M136").Value = -158055120
$ws.Range("N136").Value = -385350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 67251.67999999999
$ws.Range("I31").Value = 58231.055
$ws.Range("J31").Value = 77399.875
$ws.Range("K31").Value = 58231.055
$ws.Range("L31").Value = 77399.875
$ws.Range("M31").Value = -57936.055
$ws.Range("N31").Value = -77989.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 67251.67999999999
$ws.Range("I34").Value = 58231.055
$ws.Range("J34").Value = 77399.875
$ws.Range("K34").Value = 58231.055
$ws.Range("L34").Value = 77399.875
$ws.Range("M34").Value = -58029.055
$ws.Range("N34").Value = -77803.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13700335
$ws.Range("I58").Value = 25001588
$ws.Range("J58").Value = 1845.1818
$ws.Range("K58").Value = 25001588
$ws.Range("L58").Value = 1845.1818
$ws.Range("M58").Value = -25001385
$ws.Range("N58").Value = -2251.1818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 13700335
$ws.Range("I136").Value = 25001588
$ws.Range("J136").Value = 1845.1818
$ws.Range("K136").Value = 75004764
$ws.Range("L136").Value = 5535.5454
$ws.Range("M136").Value = -75002214
$ws.Range("N136").Value = -10635.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 727.4545000000001
$ws.Range("I5").Value = 448.73914
$ws.Range("J5").Value = 1368.5
$ws.Range("K5").Value = 1346.21742
$ws.Range("L5").Value = 4105.5
$ws.Range("M5").Value = -1234.21742
$ws.Range("N5").Value = -4329.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 767.7143
$ws.Range("I113").Value = 360
$ws.Range("K113").Value = 1080
$ws.Range("M113").Value = 1090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 855.86
$ws.Range("I131").Value = 480.9
$ws.Range("J131").Value = 949.6
$ws.Range("K131").Value = 1442.7
$ws.Range("L131").Value = 2848.8
$ws.Range("M131").Value = 3597.3
$ws.Range("N131").Value = -12928.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 727.4545000000001
$ws.Range("I135").Value = 448.73914
$ws.Range("J135").Value = 1368.5
$ws.Range("K135").Value = 4038.65226
$ws.Range("L135").Value = 12316.5
$ws.Range("M135").Value = -1503.65226
$ws.Range("N135").Value = -17386.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1002.03705
$ws.Range("I16").Value = 850.4286
$ws.Range("J16").Value = 1532.6666
$ws.Range("K16").Value = 850.4286
$ws.Range("L16").Value = 1532.6666
$ws.Range("M16").Value = -680.4286
$ws.Range("N16").Value = -1872.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1038.5883
$ws.Range("J22").Value = 1000.6
$ws.Range("L22").Value = 1000.6
$ws.Range("N22").Value = -1590.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1038.5883
$ws.Range("J27").Value = 1000.6
$ws.Range("L27").Value = 1000.6
$ws.Range("N27").Value = -1214.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 213.58333
$ws.Range("I55").Value = 221.71428
$ws.Range("J55").Value = 156.66667
$ws.Range("K55").Value = 221.71428
$ws.Range("L55").Value = 156.66667
$ws.Range("M55").Value = -48.71428
$ws.Range("N55").Value = -502.66667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2125.1875
$ws.Range("I61").Value = 2089.8
$ws.Range("K61").Value = 2089.8
$ws.Range("M61").Value = -1887.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2299.182
$ws.Range("I82").Value = 1156.6666
$ws.Range("K82").Value = 1156.6666
$ws.Range("M82").Value = -795.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2299.182
$ws.Range("I85").Value = 1156.6666
$ws.Range("K85").Value = 1156.6666
$ws.Range("M85").Value = 91.33339999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1433.7646
$ws.Range("I93").Value = 1466.9231
$ws.Range("J93").Value = 1326
$ws.Range("K93").Value = 1466.9231
$ws.Range("L93").Value = 1326
$ws.Range("M93").Value = -218.9231
$ws.Range("N93").Value = -3822

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1777.6666
$ws.Range("I100").Value = 1371.1428
$ws.Range("K100").Value = 1371.1428
$ws.Range("M100").Value = -830.1428000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2125.1875
$ws.Range("I113").Value = 2089.8
$ws.Range("K113").Value = 2089.8
$ws.Range("M113").Value = 80.19999999999982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 47300.418
$ws.Range("I132").Value = 5107.9287
$ws.Range("J132").Value = 106369.9
$ws.Range("K132").Value = 15323.7861
$ws.Range("L132").Value = 319109.7
$ws.Range("M132").Value = -12793.7861
$ws.Range("N132").Value = -324169.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 39066.145
$ws.Range("I136").Value = 22882.432
$ws.Range("J136").Value = 204140
$ws.Range("K136").Value = 68647.296
$ws.Range("L136").Value = 612420
$ws.Range("M136").Value = -66097.296
$ws.Range("N136").Value = -617520

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 883.6818
$ws.Range("I113").Value = 958.875
$ws.Range("K113").Value = 2876.625
$ws.Range("M113").Value = -706.625
